# Update the ACE AVIATION course row:
#  - C2 ("department") changes from "ACE AVIATION" to "AVIATION"
#  - R2 ("promotionValidity") is cleared (was "Promotion valid until  31th Dec 2021")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "AVIATION"
$ws.Range("R2").ClearContents()
